# Update Kosovo MSME country indicator figures with more precise values.
# These figures are stored as text in the sheet (not numbers), so the
# number format is temporarily switched to Text while writing the new
# value, then restored to the cell's original style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B13" = "57.48"
    "C13" = "0.92"
    "D14" = "62.24"
    "B16" = "98.38"
    "C16" = "1.57"
    "D16" = "99.94"
    "B20" = "16.79"
    "C20" = "26.51"
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$addr]
    $rng.Style = $origStyle
}
